# Adiciona os registros de chamada do dia 17/08/2025 e 24/08/2025
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chamada")

$attendance = @(
    @('OsValdina Francisca', 'Presente', '17/08/2025'),
    @('Paulo Henrique', 'Ausente', '17/08/2025'),
    @('João Vitor', 'Presente', '17/08/2025'),
    @('Elza Alves', 'Presente', '17/08/2025'),
    @('Antonio Patricio', 'Presente', '17/08/2025'),
    @('Gesmindo Boostel', 'Presente', '17/08/2025'),
    @('Kalahan Boostel', 'Presente', '17/08/2025'),
    @('Geciel Polegario', 'Presente', '17/08/2025'),
    @('Diana', 'Ausente', '17/08/2025'),
    @('Vanuza Nascimento', 'Presente', '17/08/2025'),
    @('Welington Nascimento', 'Ausente', '17/08/2025'),
    @('Welington Ribeiro', 'Ausente', '17/08/2025'),
    @('Jorge', 'Presente', '17/08/2025'),
    @('Gosmira', 'Presente', '17/08/2025'),
    @('Almir Rodrigues', 'Presente', '17/08/2025'),
    @('OsValdina Francisca', 'Presente', '24/08/2025'),
    @('Paulo Henrique', 'Ausente', '24/08/2025'),
    @('João Vitor', 'Ausente', '24/08/2025'),
    @('Elza Alves', 'Presente', '24/08/2025'),
    @('Antonio Patricio', 'Presente', '24/08/2025'),
    @('Gesmindo Boostel', 'Presente', '24/08/2025'),
    @('Kalahan Boostel', 'Presente', '24/08/2025'),
    @('Geciel Polegario', 'Ausente', '24/08/2025'),
    @('Diana', 'Ausente', '24/08/2025'),
    @('Vanuza Nascimento', 'Ausente', '24/08/2025'),
    @('Welington Nascimento', 'Ausente', '24/08/2025'),
    @('Welington Ribeiro', 'Ausente', '24/08/2025'),
    @('Jorge', 'Ausente', '24/08/2025'),
    @('Gosmira', 'Ausente', '24/08/2025'),
    @('Almir Rodrigues', 'Presente', '24/08/2025')
)

$startRow = 227
for ($i = 0; $i -lt $attendance.Count; $i++) {
    $row = $startRow + $i
    $record = $attendance[$i]
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
}
